$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row additions: C1 "Species", D1 "Reference" ---
$ws.Range("C1").Value = "Species"
$ws.Range("D1").Value = "Reference"

# Match the existing header formatting (bold font, border, centered/top
# aligned) used by A1:B1 by copying the format (xlPasteFormats = -4122,
# values are left untouched) from B1 onto the new header cells.
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: new Species / Reference columns, left blank ---
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""

# --- Row 3: brand new record ---
$ws.Range("A3").Value = "C:\Users\Veeraraju_elluru\Desktop\Veeraraju\Personal\ATREE\test_images\cv_ex.jpg"
$ws.Range("B3").Value = "Chihuahua, pictures"
$ws.Range("C3").Value = "sheep"
$ws.Range("D3").Value = "https://en.wikipedia.org"
